# Remove the "[ Name ]" speaker-tag runs from the outline's plan-of-presentation
# bullet points. Each tag is merged away together with exactly one trailing
# space that used to separate it from the preceding text (Find/Replace keeps
# the formatting of the first run in the match, which is exactly the
# desired result here).

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Plan of presentation  [ Olga ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Plan of presentation ", 2) | Out-Null

$d.Content.Find.Execute(
    "History  [ Steve ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "History ", 2) | Out-Null

$d.Content.Find.Execute(
    "The problem statement and real life examples. [Konstantin]", $false, $false, $false, $false, $false,
    $true, 1, $false, "The problem statement and real life examples.", 2) | Out-Null

$d.Content.Find.Execute(
    "Audience will pick up nodes with weight  [ Steve ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Audience will pick up nodes with weight ", 2) | Out-Null

$d.Content.Find.Execute(
    "Terminology / Initialization / Algorithm    [ Steve ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Terminology / Initialization / Algorithm   ", 2) | Out-Null

$d.Content.Find.Execute(
    "Solve the problem step by step with graph [ Steve ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Solve the problem step by step with graph", 2) | Out-Null

# The "Show code and Graph" bullet only loses its trailing "[ Olga ]" run;
# the text before it (including the lone space run) stays untouched. By now
# the other "[ Olga ]" tag has already been consumed above, so this only
# matches the remaining one.
$d.Content.Find.Execute(
    "[ Olga ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 2) | Out-Null

$d.Content.Find.Execute(
    "Analyze runtime  [ Konstantin ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Analyze runtime ", 2) | Out-Null

$d.Content.Find.Execute(
    "Improvements (Implementation using heap/priority queue + use adjacency list)  [ Konstantin ]", $false, $false, $false, $false, $false,
    $true, 1, $false, "Improvements (Implementation using heap/priority queue + use adjacency list) ", 2) | Out-Null

Write-Output "Removed speaker-tag runs from outline bullets"
